$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), copying the header style from G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values per row (2..22): 1 for the two big outings (rows 8 & 13), 0 otherwise
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
}

foreach ($row in 2..22) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
